$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Carrier" column header and its row-2 value.
$ws.Range("G1").Value = "Carrier"
$ws.Range("G2").Value = "Echo,Team Worldwide"

# Update the selection to match the authored change.
$ws.Range("E2").Select()
